$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '310.10'

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-2.58%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '37.67'

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-5.07%'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.105'

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.70%'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07858'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-4.30%'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.966'

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-3.66%'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '4.369'

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.98%'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '8.304'

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.02%'

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-5.40%'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9281'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.38%'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1347'

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-4.44%'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1990'

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.19%'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.08940'

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-1.53%'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03466'

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.28%'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09697'

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-1.10%'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001388'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.40%'

$ws.Range("B17").Value = 'CoinExToken'

$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04338'

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.34%'

$ws.Range("B18").Value = 'TigerCash'

$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.005967'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-3.81%'

$ws.Range("B19").Value = 'UpBots'

$ws.Range("C19").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.007506'

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1,778.08%'

$ws.Range("B20").Value = 'LEO'

$ws.Range("C20").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.587'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-2.66%'

$ws.Range("B21").Value = 'BitpandaEcosystemToken'

$ws.Range("C21").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.3466'

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.26%'

$ws.Range("B22").Value = 'ProBitToken'

$ws.Range("C22").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.1295'

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.02%'

$ws.Range("B23").Value = 'MCDex'

$ws.Range("C23").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.002'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2.17%'

$ws.Range("B24").Value = 'ZBToken'

$ws.Range("C24").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.2511'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '2.55%'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001226'

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.08%'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004545'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-4.71%'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001351'

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '4.00%'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02299'

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '3.84%'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05065'

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-3.17%'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007473'

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.46%'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009829'

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '1.71%'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1358'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-1.50%'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.001980'

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-6.28%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.008779'

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-10.89%'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006828'

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '3.48%'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000750'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.14%'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003001'

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '8.62%'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001301'

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '8.38%'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002101'

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.14%'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002001'

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.14%'
